$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the header row (row 1) from English emoji descriptors to
# Vietnamese words, per the "remove stopwords / remove _" cleanup.
$ws.Range("A1").Value = "lè lưỡi"
$ws.Range("B1").Value = "cười lăn lộn"
$ws.Range("C1").Value = "cười mĩm"
$ws.Range("D1").Value = "cười"
$ws.Range("E1").Value = "cười ra nước mắt"
$ws.Range("F1").Value = "cảm thấy buồn"
$ws.Range("G1").Value = "muốn khóc"
$ws.Range("H1").Value = "mặt nhăn"
$ws.Range("I1").Value = "rối rắm"
$ws.Range("J1").Value = "cạn lời"
$ws.Range("K1").Value = "cười híp mắt"

# Update the active selection to B2, as reflected in the saved sheet view.
$ws.Range("B2").Select()
